# Apply updated cryptocurrency price/volume data to the worksheet.
# Each target cell stores plain text (t="inlineStr" in the original file),
# so we force a text number format before assigning the value to stop Excel
# from reinterpreting numeric-looking strings (e.g. "8.90") as numbers, then
# restore the default "Normal" style so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "44.349.92"
Set-TextValue "E2" "  +2.42%  "
Set-TextValue "D3" "2.364.18"
Set-TextValue "E3" "  +0.43%  "
Set-TextValue "E4" "  +0.08%  "
Set-TextValue "D5" "0.677"
Set-TextValue "E5" "  +4.10%  "
Set-TextValue "D6" "238.76"
Set-TextValue "E6" "  +3.06%  "
Set-TextValue "D7" "73.39"
Set-TextValue "E7" "  +10.25%  "
Set-TextValue "E8" "  +0.02%  "
Set-TextValue "D9" "0.547"
Set-TextValue "E9" "  +19.58%  "
Set-TextValue "D10" "0.103"
Set-TextValue "E10" "  +9.29%  "
Set-TextValue "D11" "29.47"
Set-TextValue "E11" "  +10.55%  "
Set-TextValue "D12" "0.107"
Set-TextValue "E12" "  +2.61%  "
Set-TextValue "D13" "2.717.55"
Set-TextValue "E13" "  +0.48%  "
Set-TextValue "D14" "16.84"
Set-TextValue "E14" "  +9.51%  "
Set-TextValue "D15" "6.71"
Set-TextValue "E15" "  +7.23%  "
Set-TextValue "D16" "0.906"
Set-TextValue "E16" "  +8.47%  "
Set-TextValue "D17" "2.364.02"
Set-TextValue "E17" "  +0.37%  "
Set-TextValue "D18" "44.305.05"
Set-TextValue "E18" "  +2.28%  "
Set-TextValue "D19" "0.0000104"
Set-TextValue "E19" "  +6.13%  "
Set-TextValue "D20" "77.94"
Set-TextValue "E20" "  +5.92%  "
Set-TextValue "D21" "6.45"
Set-TextValue "E21" "  +3.91%  "
Set-TextValue "D22" "255.66"
Set-TextValue "E22" "  +3.13%  "
Set-TextValue "E23" "  -0.14%  "
Set-TextValue "E24" "  -4.01%  "
Set-TextValue "D25" "2.52"
Set-TextValue "E25" "  +3.36%  "
Set-TextValue "D26" "10.51"
Set-TextValue "E26" "  +5.75%  "
Set-TextValue "D27" "2.24"
Set-TextValue "E27" "  -1.22%  "
Set-TextValue "D28" "22.53"
Set-TextValue "E28" "  +0.87%  "
Set-TextValue "D29" "173.18"
Set-TextValue "E29" "  -1.03%  "
Set-TextValue "E30" "  +5.08%  "
Set-TextValue "E31" "  +3.53%  "
Set-TextValue "E32" "  +5.23%  "
Set-TextValue "D33" "0.0741"
Set-TextValue "E33" "  +6.87%  "
Set-TextValue "D34" "5.21"
Set-TextValue "E34" "  +4.87%  "
Set-TextValue "D35" "5.22"
Set-TextValue "E35" "  +4.33%  "
Set-TextValue "D36" "3.93"
Set-TextValue "E36" "  +8.99%  "
Set-TextValue "D37" "2.44"
Set-TextValue "E37" "  -2.68%  "
Set-TextValue "E38" "  +0.78%  "
Set-TextValue "D39" "0.0271"
Set-TextValue "E39" "  +7.48%  "
Set-TextValue "D40" "19.74"
Set-TextValue "E40" "  +10.92%  "
Set-TextValue "E41" "  +0.03%  "
Set-TextValue "D42" "8.90"
Set-TextValue "E42" "  -0.24%  "
Set-TextValue "D43" "1.26"
Set-TextValue "E43" "  +3.72%  "
Set-TextValue "D44" "0.0980"
Set-TextValue "E44" "  +3.57%  "
Set-TextValue "E45" "  +0.45%  "
Set-TextValue "D46" "4.49"
Set-TextValue "E46" "  +3.36%  "
Set-TextValue "D47" "98.60"
Set-TextValue "E47" "  -0.07%  "
Set-TextValue "E48" "  +12.58%  "
Set-TextValue "E49" "  +5.17%  "
Set-TextValue "D50" "1.442.77"
Set-TextValue "E50" "  +0.34%  "
Set-TextValue "B51" "HuobiToken"
Set-TextValue "C51" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D51" "2.78"
Set-TextValue "E51" "  +1.53%  "
